# ---------------------------------------------------------------------------
# Adds a new "2022-Q3" sheet (with fund-holdings detail) right after "总计"
# and prepends a matching summary row on the "总计" sheet.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new row 2 for "2022-Q3" and push the
#    rest of the quarters down by one row.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()
$summary.Range("B2:D2").ClearFormats()

# Re-use the index-column style (bold + border) from the row below, which
# still carries the original formatting for column A.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q3"
$summary.Cells.Item(2,3).Value = 30
$summary.Cells.Item(2,4).Value = 2.39

# The index column (A) is a simple 0-based row counter; renumber the rows
# that got pushed down by the insert.
$summary.Cells.Item(3,1).Value = 1
$summary.Cells.Item(4,1).Value = 2
$summary.Cells.Item(5,1).Value = 3
$summary.Cells.Item(6,1).Value = 4
$summary.Cells.Item(7,1).Value = 5

# ---------------------------------------------------------------------------
# 2) New "2022-Q3" worksheet: duplicate the "2022-Q2" sheet (so headers,
#    column styling and borders match the other quarterly sheets) and place
#    it immediately after "总计" / before "2022-Q2".
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2, $null)
$ws = $wb.Worksheets.Item("2022-Q2 (2)")
$ws.Name = "2022-Q3"

# Extend the bold/bordered index-column style (column A) down through row 31
# so every data row matches the look of the other quarterly sheets.
$ws.Range("A2").Copy()
$ws.Range("A5:A31").PasteSpecial(-4122)

# Clear the 3 sample rows that came over from the "2022-Q2" copy so they
# don't leave stray formatting/content behind (they get fully repopulated
# below).
$ws.Range("B2:H4").ClearContents()

$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "'013141"
$ws.Cells.Item(2,3).Value = '中信保诚弘远混合A'
$ws.Cells.Item(2,4).Value = "'19.17"
$ws.Cells.Item(2,5).Value = "'79.40"
$ws.Cells.Item(2,6).Value = "'3.01"
$ws.Cells.Item(2,7).Value = "'0.5770"
$ws.Cells.Item(2,8).Value = 8
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = "'001210"
$ws.Cells.Item(3,3).Value = '天弘互联网灵活配置混合A'
$ws.Cells.Item(3,4).Value = "'7.93"
$ws.Cells.Item(3,5).Value = "'92.94"
$ws.Cells.Item(3,6).Value = "'3.67"
$ws.Cells.Item(3,7).Value = "'0.2910"
$ws.Cells.Item(3,8).Value = 10
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = "'050010"
$ws.Cells.Item(4,3).Value = '博时特许价值混合'
$ws.Cells.Item(4,4).Value = "'6.01"
$ws.Cells.Item(4,5).Value = "'92.93"
$ws.Cells.Item(4,6).Value = "'3.49"
$ws.Cells.Item(4,7).Value = "'0.2097"
$ws.Cells.Item(4,8).Value = 10
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = "'166011"
$ws.Cells.Item(5,3).Value = '中欧盛世成长混合（LOF）A'
$ws.Cells.Item(5,4).Value = "'5.21"
$ws.Cells.Item(5,5).Value = "'81.90"
$ws.Cells.Item(5,6).Value = "'3.57"
$ws.Cells.Item(5,7).Value = "'0.1860"
$ws.Cells.Item(5,8).Value = 6
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = "'000462"
$ws.Cells.Item(6,3).Value = '农银主题轮动灵活配置混合'
$ws.Cells.Item(6,4).Value = "'4.69"
$ws.Cells.Item(6,5).Value = "'91.45"
$ws.Cells.Item(6,6).Value = "'3.20"
$ws.Cells.Item(6,7).Value = "'0.1501"
$ws.Cells.Item(6,8).Value = 8
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = "'001306"
$ws.Cells.Item(7,3).Value = '中欧永裕混合A'
$ws.Cells.Item(7,4).Value = "'4.05"
$ws.Cells.Item(7,5).Value = "'81.72"
$ws.Cells.Item(7,6).Value = "'3.57"
$ws.Cells.Item(7,7).Value = "'0.1446"
$ws.Cells.Item(7,8).Value = 6
$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = "'006058"
$ws.Cells.Item(8,3).Value = '民生加银新兴成长混合'
$ws.Cells.Item(8,4).Value = "'3.98"
$ws.Cells.Item(8,5).Value = "'87.53"
$ws.Cells.Item(8,6).Value = "'3.22"
$ws.Cells.Item(8,7).Value = "'0.1282"
$ws.Cells.Item(8,8).Value = 10
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = "'001463"
$ws.Cells.Item(9,3).Value = '光大保德信一带一路战略主题混合'
$ws.Cells.Item(9,4).Value = "'1.57"
$ws.Cells.Item(9,5).Value = "'87.51"
$ws.Cells.Item(9,6).Value = "'5.94"
$ws.Cells.Item(9,7).Value = "'0.0933"
$ws.Cells.Item(9,8).Value = 4
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = "'501200"
$ws.Cells.Item(10,3).Value = '民生加银科技创新 3 年封闭运作灵活配置混合'
$ws.Cells.Item(10,4).Value = "'2.62"
$ws.Cells.Item(10,5).Value = "'87.51"
$ws.Cells.Item(10,6).Value = "'3.29"
$ws.Cells.Item(10,7).Value = "'0.0862"
$ws.Cells.Item(10,8).Value = 9
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = "'000039"
$ws.Cells.Item(11,3).Value = '农银高增长混合'
$ws.Cells.Item(11,4).Value = "'2.51"
$ws.Cells.Item(11,5).Value = "'90.96"
$ws.Cells.Item(11,6).Value = "'3.19"
$ws.Cells.Item(11,7).Value = "'0.0801"
$ws.Cells.Item(11,8).Value = 8
$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = "'002707"
$ws.Cells.Item(12,3).Value = '摩根士丹利华鑫科技领先灵活配置混合A'
$ws.Cells.Item(12,4).Value = "'1.75"
$ws.Cells.Item(12,5).Value = "'94.13"
$ws.Cells.Item(12,6).Value = "'4.02"
$ws.Cells.Item(12,7).Value = "'0.0704"
$ws.Cells.Item(12,8).Value = 5
$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = "'012924"
$ws.Cells.Item(13,3).Value = '华夏新时代灵活配置混合（QDII）美元现汇'
$ws.Cells.Item(13,4).Value = "'2.11"
$ws.Cells.Item(13,5).Value = "'73.45"
$ws.Cells.Item(13,6).Value = "'3.04"
$ws.Cells.Item(13,7).Value = "'0.0641"
$ws.Cells.Item(13,8).Value = 5
$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(14,2).Value = "'012925"
$ws.Cells.Item(14,3).Value = '华夏新时代灵活配置混合（QDII）美元现钞'
$ws.Cells.Item(14,4).Value = "'2.11"
$ws.Cells.Item(14,5).Value = "'73.45"
$ws.Cells.Item(14,6).Value = "'3.04"
$ws.Cells.Item(14,7).Value = "'0.0641"
$ws.Cells.Item(14,8).Value = 5
$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(15,2).Value = "'012259"
$ws.Cells.Item(15,3).Value = '天弘鑫悦成长混合C'
$ws.Cells.Item(15,4).Value = "'1.19"
$ws.Cells.Item(15,5).Value = "'93.52"
$ws.Cells.Item(15,6).Value = "'4.76"
$ws.Cells.Item(15,7).Value = "'0.0566"
$ws.Cells.Item(15,8).Value = 8
$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(16,2).Value = "'001983"
$ws.Cells.Item(16,3).Value = '中邮低碳经济灵活配置混合'
$ws.Cells.Item(16,4).Value = "'0.50"
$ws.Cells.Item(16,5).Value = "'91.21"
$ws.Cells.Item(16,6).Value = "'6.14"
$ws.Cells.Item(16,7).Value = "'0.0307"
$ws.Cells.Item(16,8).Value = 6
$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(17,2).Value = "'007713"
$ws.Cells.Item(17,3).Value = '华富科技动能混合'
$ws.Cells.Item(17,4).Value = "'0.60"
$ws.Cells.Item(17,5).Value = "'89.23"
$ws.Cells.Item(17,6).Value = "'4.29"
$ws.Cells.Item(17,7).Value = "'0.0257"
$ws.Cells.Item(17,8).Value = 9
$ws.Cells.Item(18,1).Value = 16
$ws.Cells.Item(18,2).Value = "'004233"
$ws.Cells.Item(18,3).Value = '中欧盛世成长混合（LOF）C'
$ws.Cells.Item(18,4).Value = "'0.64"
$ws.Cells.Item(18,5).Value = "'81.90"
$ws.Cells.Item(18,6).Value = "'3.57"
$ws.Cells.Item(18,7).Value = "'0.0228"
$ws.Cells.Item(18,8).Value = 6
$ws.Cells.Item(19,1).Value = 17
$ws.Cells.Item(19,2).Value = "'015005"
$ws.Cells.Item(19,3).Value = '中邮能源革新混合C'
$ws.Cells.Item(19,4).Value = "'0.40"
$ws.Cells.Item(19,5).Value = "'91.38"
$ws.Cells.Item(19,6).Value = "'4.91"
$ws.Cells.Item(19,7).Value = "'0.0196"
$ws.Cells.Item(19,8).Value = 7
$ws.Cells.Item(20,1).Value = 18
$ws.Cells.Item(20,2).Value = "'008998"
$ws.Cells.Item(20,3).Value = '同泰竞争优势混合C'
$ws.Cells.Item(20,4).Value = "'0.59"
$ws.Cells.Item(20,5).Value = "'93.98"
$ws.Cells.Item(20,6).Value = "'3.30"
$ws.Cells.Item(20,7).Value = "'0.0195"
$ws.Cells.Item(20,8).Value = 3
$ws.Cells.Item(21,1).Value = 19
$ws.Cells.Item(21,2).Value = "'001888"
$ws.Cells.Item(21,3).Value = '中欧盛世成长混合（LOF）E'
$ws.Cells.Item(21,4).Value = "'0.49"
$ws.Cells.Item(21,5).Value = "'81.90"
$ws.Cells.Item(21,6).Value = "'3.57"
$ws.Cells.Item(21,7).Value = "'0.0175"
$ws.Cells.Item(21,8).Value = 6
$ws.Cells.Item(22,1).Value = 20
$ws.Cells.Item(22,2).Value = "'012258"
$ws.Cells.Item(22,3).Value = '天弘鑫悦成长混合A'
$ws.Cells.Item(22,4).Value = "'0.28"
$ws.Cells.Item(22,5).Value = "'93.52"
$ws.Cells.Item(22,6).Value = "'4.76"
$ws.Cells.Item(22,7).Value = "'0.0133"
$ws.Cells.Item(22,8).Value = 8
$ws.Cells.Item(23,1).Value = 21
$ws.Cells.Item(23,2).Value = "'001307"
$ws.Cells.Item(23,3).Value = '中欧永裕混合C'
$ws.Cells.Item(23,4).Value = "'0.33"
$ws.Cells.Item(23,5).Value = "'81.72"
$ws.Cells.Item(23,6).Value = "'3.57"
$ws.Cells.Item(23,7).Value = "'0.0118"
$ws.Cells.Item(23,8).Value = 6
$ws.Cells.Item(24,1).Value = 22
$ws.Cells.Item(24,2).Value = "'008997"
$ws.Cells.Item(24,3).Value = '同泰竞争优势混合A'
$ws.Cells.Item(24,4).Value = "'0.26"
$ws.Cells.Item(24,5).Value = "'93.98"
$ws.Cells.Item(24,6).Value = "'3.30"
$ws.Cells.Item(24,7).Value = "'0.0086"
$ws.Cells.Item(24,8).Value = 3
$ws.Cells.Item(25,1).Value = 23
$ws.Cells.Item(25,2).Value = "'015004"
$ws.Cells.Item(25,3).Value = '中邮能源革新混合A'
$ws.Cells.Item(25,4).Value = "'0.12"
$ws.Cells.Item(25,5).Value = "'91.38"
$ws.Cells.Item(25,6).Value = "'4.91"
$ws.Cells.Item(25,7).Value = "'0.0059"
$ws.Cells.Item(25,8).Value = 7
$ws.Cells.Item(26,1).Value = 24
$ws.Cells.Item(26,2).Value = "'015461"
$ws.Cells.Item(26,3).Value = '天弘互联网灵活配置混合C'
$ws.Cells.Item(26,4).Value = "'0.13"
$ws.Cells.Item(26,5).Value = "'92.94"
$ws.Cells.Item(26,6).Value = "'3.67"
$ws.Cells.Item(26,7).Value = "'0.0048"
$ws.Cells.Item(26,8).Value = 10
$ws.Cells.Item(27,1).Value = 25
$ws.Cells.Item(27,2).Value = "'005281"
$ws.Cells.Item(27,3).Value = '中科沃土转型升级灵活配置混合'
$ws.Cells.Item(27,4).Value = "'0.10"
$ws.Cells.Item(27,5).Value = "'57.70"
$ws.Cells.Item(27,6).Value = "'4.00"
$ws.Cells.Item(27,7).Value = "'0.0040"
$ws.Cells.Item(27,8).Value = 1
$ws.Cells.Item(28,1).Value = 26
$ws.Cells.Item(28,2).Value = "'014871"
$ws.Cells.Item(28,3).Value = '摩根士丹利华鑫科技领先灵活配置混合C'
$ws.Cells.Item(28,4).Value = "'0.08"
$ws.Cells.Item(28,5).Value = "'94.13"
$ws.Cells.Item(28,6).Value = "'4.02"
$ws.Cells.Item(28,7).Value = "'0.0032"
$ws.Cells.Item(28,8).Value = 5
$ws.Cells.Item(29,1).Value = 27
$ws.Cells.Item(29,2).Value = "'004522"
$ws.Cells.Item(29,3).Value = '安信工业4.0主题沪港深精选灵活配置混合C'
$ws.Cells.Item(29,4).Value = "'0.07"
$ws.Cells.Item(29,5).Value = "'75.24"
$ws.Cells.Item(29,6).Value = "'3.10"
$ws.Cells.Item(29,7).Value = "'0.0022"
$ws.Cells.Item(29,8).Value = 9
$ws.Cells.Item(30,1).Value = 28
$ws.Cells.Item(30,2).Value = "'004521"
$ws.Cells.Item(30,3).Value = '安信工业4.0主题沪港深精选灵活配置混合A'
$ws.Cells.Item(30,4).Value = "'0.03"
$ws.Cells.Item(30,5).Value = "'75.24"
$ws.Cells.Item(30,6).Value = "'3.10"
$ws.Cells.Item(30,7).Value = "'0.0009"
$ws.Cells.Item(30,8).Value = 9
$ws.Cells.Item(31,1).Value = 29
$ws.Cells.Item(31,2).Value = "'015936"
$ws.Cells.Item(31,3).Value = '中信保诚弘远混合C'
$ws.Cells.Item(31,4).Value = "'0.02"
$ws.Cells.Item(31,5).Value = "'79.40"
$ws.Cells.Item(31,6).Value = "'3.01"
$ws.Cells.Item(31,7).Value = "'0.0006"
$ws.Cells.Item(31,8).Value = 8
